# Update excess mortality CBS model
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct existing data: week 53 (2020) and week 1 (2021) observed counts ---
$ws.Range("H45").Value = 4077
$ws.Range("H46").Value = 4118

# --- Add newly observed weeks for 2021 (weeks 2, 3 and 4) ---
$ws.Range("F47").Value = 2021
$ws.Range("G47").Value = 2
$ws.Range("H47").Value = 3825
$ws.Range("I47").Value = 3343
$ws.Range("J47").Formula = "=H47-I47"

$ws.Range("F48").Value = 2021
$ws.Range("G48").Value = 3
$ws.Range("H48").Value = 3823
$ws.Range("I48").Value = 3376
$ws.Range("J48").Formula = "=H48-I48"

$ws.Range("F49").Value = 2021
$ws.Range("G49").Value = 4
$ws.Range("H49").Value = 3758
$ws.Range("I49").Value = 3425
$ws.Range("J49").Formula = "=H49-I49"

# --- Drop the old "Som week 11 tot en met 19" summary row and replace it
#     with a blank (but still numerically formatted) row 51 ---
$ws.Rows("50:50").Delete()
$ws.Range("H51").NumberFormat = "0"
$ws.Range("I51").NumberFormat = "0"
$ws.Range("J51").NumberFormat = "0"

# --- Restore view state ---
$ws.Range("A15").Select()
$ws.Range("H37").Select()
